$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.719.09'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.154.73'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.00'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.78'
$ws.Range("E6").Value = '  +1.02%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  +14.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.428'
$ws.Range("E10").Value = '  +5.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.110'
$ws.Range("E11").Value = '  +2.54%  '
$ws.Range("E12").Value = '  +2.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.698.69'
$ws.Range("E13").Value = '  +1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.82'
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("E15").Value = '  +4.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.778.31'
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.21'
$ws.Range("E17").Value = '  +3.95%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.146.51'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.00'
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.19'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.43'
$ws.Range("E21").Value = '  +5.57%  '
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.93'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.519'
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.168'
$ws.Range("E26").Value = '  +1.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.20'
$ws.Range("E28").Value = '  +13.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0863'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.10'
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.95'
$ws.Range("E32").Value = '  +3.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.20'
$ws.Range("E33").Value = '  +4.95%  '
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.20'
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.28'
$ws.Range("E36").Value = '  +4.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.34'
$ws.Range("E37").Value = '  +7.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.30'
$ws.Range("E38").Value = '  -3.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.650.97'
$ws.Range("E39").Value = '  +10.77%  '
$ws.Range("E40").Value = '  +1.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0681'
$ws.Range("E41").Value = '  +2.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.23'
$ws.Range("E42").Value = '  +6.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.77'
$ws.Range("E43").Value = '  +3.60%  '
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0286'
$ws.Range("E45").Value = '  +7.85%  '
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.196.42'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E48").Value = '  +13.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.985'
$ws.Range("E49").Value = '  +2.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.19'
$ws.Range("E50").Value = '  +2.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.10'
$ws.Range("E51").Value = '  +2.11%  '
